$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-storage (matches source data: numeric-looking IDs/dates/codes are stored as text)
$ws.Range("A105:H109").NumberFormat = "@"
$ws.Range("J105:L109").NumberFormat = "@"
$ws.Range("O105:R109").NumberFormat = "@"

# Row 105
$ws.Range("A105").Value = '5931'
$ws.Range("B105").Value = '11/17/2025'
$ws.Range("C105").Value = 'Paso 280'
$ws.Range("D105").Value = '3'
$ws.Range("E105").Value = 'Pendiente ADM'
$ws.Range("F105").Value = 'PEBCOM'
$ws.Range("G105").Value = 'Pendiente'
$ws.Range("H105").Value = 'base picada'
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = 'Cambio'
$ws.Range("K105").Value = 'Sin equipos'
$ws.Range("L105").Value = 'Pasante'
$ws.Range("M105").Value = -58.403566
$ws.Range("N105").Value = -34.606691
$ws.Range("O105").Value = 'Almagro'
$ws.Range("P105").Value = 'Capital Sur'
$ws.Range("Q105").Value = 'CLI-D'
$ws.Range("R105").Value = 'Fuera de Poligono OVL'

# Row 106
$ws.Range("A106").Value = '7581'
$ws.Range("B106").Value = '11/14/2025'
$ws.Range("C106").Value = 'SUAREZ 951'
$ws.Range("D106").Value = '4'
$ws.Range("E106").Value = 'Pendiente ADM'
$ws.Range("F106").Value = 'PEBCOM'
$ws.Range("G106").Value = 'Pendiente'
$ws.Range("H106").Value = 'corroida'
$ws.Range("I106").Value = 1
$ws.Range("J106").Value = 'Cambio'
$ws.Range("K106").Value = 'Sin equipos'
$ws.Range("L106").Value = 'Pasante'
$ws.Range("M106").Value = -58.36561
$ws.Range("N106").Value = -34.637763
$ws.Range("O106").Value = 'San Telmo'
$ws.Range("P106").Value = 'Capital Sur'
$ws.Range("Q106").Value = 'CON-D'
$ws.Range("R106").Value = 'Fuera de Poligono OVL'

# Row 107
$ws.Range("A107").Value = '7845'
$ws.Range("B107").Value = '11/14/2025'
$ws.Range("C107").Value = 'VELEZ SARSFIELD AV. 10'
$ws.Range("D107").Value = '4'
$ws.Range("E107").Value = 'Pendiente ADM'
$ws.Range("F107").Value = 'PEBCOM'
$ws.Range("G107").Value = 'Pendiente'
$ws.Range("H107").Value = 'columna inclinada'
$ws.Range("I107").Value = 1
$ws.Range("J107").Value = 'Cambio'
$ws.Range("K107").Value = 'Sin equipos'
$ws.Range("L107").Value = 'Pasante'
$ws.Range("M107").Value = -58.390341
$ws.Range("N107").Value = -34.634311
$ws.Range("O107").Value = 'San Telmo'
$ws.Range("P107").Value = 'Capital Sur'
$ws.Range("Q107").Value = 'CON-K'
$ws.Range("R107").Value = 'Fuera de Poligono OVL'

# Row 108
$ws.Range("A108").Value = '7846'
$ws.Range("B108").Value = '11/14/2025'
$ws.Range("C108").Value = 'DIAZ VELEZ AV. 3485'
$ws.Range("D108").Value = '5'
$ws.Range("E108").Value = 'Pendiente ADM'
$ws.Range("F108").Value = 'PEBCOM'
$ws.Range("G108").Value = 'Pendiente'
$ws.Range("H108").Value = 'columna inclinada'
$ws.Range("I108").Value = 1
$ws.Range("J108").Value = 'Cambio'
$ws.Range("K108").Value = 'Sin equipos'
$ws.Range("L108").Value = 'Pasante'
$ws.Range("M108").Value = -58.415838
$ws.Range("N108").Value = -34.608469
$ws.Range("O108").Value = 'Almagro'
$ws.Range("P108").Value = 'Capital Sur'
$ws.Range("Q108").Value = 'CLI-J'
$ws.Range("R108").Value = 'Fuera de Poligono OVL'

# Row 109
$ws.Range("A109").Value = '7850'
$ws.Range("B109").Value = '11/14/2025'
$ws.Range("C109").Value = 'ARCAMENDIA 793'
$ws.Range("D109").Value = '4'
$ws.Range("E109").Value = 'Pendiente ADM'
$ws.Range("F109").Value = 'PEBCOM'
$ws.Range("G109").Value = 'Pendiente'
$ws.Range("H109").Value = 'columna inclinada chocada doblada'
$ws.Range("I109").Value = 1
$ws.Range("J109").Value = 'Cambio'
$ws.Range("K109").Value = 'Sin equipos'
$ws.Range("L109").Value = 'Pasante'
$ws.Range("M109").Value = -58.379965
$ws.Range("N109").Value = -34.641243
$ws.Range("O109").Value = 'San Telmo'
$ws.Range("P109").Value = 'Capital Sur'
$ws.Range("Q109").Value = 'CON-H'
$ws.Range("R109").Value = 'Fuera de Poligono OVL'

# Reset style so no explicit style index is left on the new cells (matches source: unstyled data rows)
$ws.Range("A105:R109").Style = "Normal"
